# Automatic update of files.
#
# Rotates the data of rows 2-4 (row2 -> row3, row3 -> row4, row4 -> row2),
# swaps the data of rows 13/14, and fixes a typo in P8.
# All other cells (row-independent columns, other rows) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51  # column AY

function Get-RowSnapshot($rowIndex) {
    $vals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($rowIndex, $c).Value2
    }
    return $vals
}

$textColumnIndexes = @(9)  # column I ("Antal") always holds text, even digits-only values

function Set-CellSafe($targetCell, $val, $colIndex) {
    # Writing a plain "yyyy-MM-dd" string back through .Value lets Excel's
    # COM layer auto-convert it into a real date; the source file stores
    # these as literal text, so force text with a leading apostrophe when
    # the value looks like an ISO date.
    if ($val -is [string] -and $val -match '^\d{4}-\d{2}-\d{2}$') {
        $targetCell.Value = "'" + $val
    } elseif (($textColumnIndexes -contains $colIndex) -and ($val -ne $null) -and ($val -ne '')) {
        # Column I is always stored as text in the source file, but reading
        # a digits-only cell back through .Value2 already coerces it to a
        # number, so re-force it to text on write.
        $targetCell.Value = "'" + [string]$val
    } else {
        $targetCell.Value = $val
    }
}

function Set-RowFromSnapshot($rowIndex, $vals) {
    for ($c = 1; $c -le $lastCol; $c++) {
        Set-CellSafe $ws.Cells.Item($rowIndex, $c) $vals[$c - 1] $c
    }
}

# --- snapshot the rows involved before overwriting anything ---
$snapRow2 = Get-RowSnapshot 2
$snapRow3 = Get-RowSnapshot 3
$snapRow4 = Get-RowSnapshot 4
$snapRow13 = Get-RowSnapshot 13
$snapRow14 = Get-RowSnapshot 14

# --- apply the rotation: row2 <- row4, row3 <- row2, row4 <- row3 ---
Set-RowFromSnapshot 2 $snapRow4
Set-RowFromSnapshot 3 $snapRow2
Set-RowFromSnapshot 4 $snapRow3

# --- apply the swap: row13 <- row14, row14 <- row13 ---
Set-RowFromSnapshot 13 $snapRow14
Set-RowFromSnapshot 14 $snapRow13

# --- fix the typo in P8 ---
$ws.Range("P8").Value = "Värsångsflon, Jmt"
